$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.920.77'
$ws.Range('E2').Value = '  +1.90%  '

$ws.Range('D3').Value = '2.492.92'
$ws.Range('E3').Value = '  +0.14%  '

$ws.Range('E4').Value = '  -0.11%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '590.74'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.73%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '174.81'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.80%  '

$ws.Range('E7').Value = '  -0.09%  '

$ws.Range('E8').Value = '  -0.27%  '

$ws.Range('D9').Value = '2.493.07'
$ws.Range('E9').Value = '  +0.17%  '

$ws.Range('E10').Value = '  +7.48%  '

$ws.Range('E11').Value = '  -0.65%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.96'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.77%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.335'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.26%  '

$ws.Range('D14').Value = '2.960.16'
$ws.Range('E14').Value = '  +0.34%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '25.53'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.58%  '

$ws.Range('D16').Value = '68.742.89'
$ws.Range('E16').Value = '  +1.72%  '

$ws.Range('E17').Value = '  +0.91%  '

$ws.Range('D18').Value = '2.492.13'
$ws.Range('E18').Value = '  +2.04%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '358.77'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.27%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.49'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.57%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.82'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.35%  '

$ws.Range('E22').Value = '  -2.04%  '

$ws.Range('E23').Value = '  -0.03%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '69.80'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.11%  '

$ws.Range('E25').Value = '  -3.16%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.88'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.68%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.65'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -6.79%  '

$ws.Range('D28').Value = '2.625.78'
$ws.Range('E28').Value = '  +0.19%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.984'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.21%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '505.43'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.33%  '

$ws.Range('D31').Value = '0.0₃0873'
$ws.Range('E31').Value = '  -3.21%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.66'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.74%  '

$ws.Range('E33').Value = '  -0.49%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.21'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.11%  '

$ws.Range('E35').Value = '  -0.02%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '161.50'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.85%  '

$ws.Range('E37').Value = '  -3.57%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '18.51'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.03%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.66'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.05%  '

$ws.Range('E40').Value = '  -0.05%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.29'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.19%  '

$ws.Range('E42').Value = '  -3.07%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.72'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.39%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.316'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.77%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.27'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -5.54%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '148.79'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.65%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.52'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.32%  '

$ws.Range('E48').Value = '  -1.44%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0731'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.64%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.54'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.61%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.572'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.30%  '
